$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (columns D, L, M, N, O, P, S)
# Column order in each row array: D, L, M, N, O, P, S
$rows = @{
    2  = @("44841", "Primera",  60,  23000, 24000, 23500, 2350)
    3  = @("44848", "Especial", 60,  24000, 25000, 24500, 2450)
    4  = @("44848", "Primera",  120, 21000, 22000, 21500, 2150)
    5  = @("44447", "Primera",  60,  21000, 22000, 21500, 2150)
    6  = @("44461", "Especial", 60,  31000, 32000, 31500, 3150)
    7  = @("44461", "Primera",  30,  30000, 30000, 30000, 3000)
    8  = @("44460", "Especial", 60,  31000, 32000, 31500, 3150)
    9  = @("44460", "Primera",  30,  30000, 30000, 30000, 3000)
    10 = @("44874", "Especial", 30,  25000, 25000, 25000, 2500)
    11 = @("44874", "Primera",  80,  23000, 24000, 23500, 2350)
    13 = @("44446", "Primera",  60,  21000, 22000, 21500, 2150)
    14 = @("44839", "Primera",  120, 25000, 26000, 25500, 2550)
    16 = @("44868", "Especial", 60,  26000, 26000, 26000, 2600)
    17 = @("44487", "Primera",  30,  23000, 24000, 23500, 2350)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = [double]$vals[0]
    $ws.Range("L$r").Value = $vals[1]
    $ws.Range("M$r").Value = $vals[2]
    $ws.Range("N$r").Value = $vals[3]
    $ws.Range("O$r").Value = $vals[4]
    $ws.Range("P$r").Value = $vals[5]
    $ws.Range("S$r").Value = $vals[6]
}
